$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for rows 92-97, columns D-L (mean_TP, mean_FP, mean_TN, mean_FN, Recall, Precision, f1, f0.5, accuracy)
$data = @{
    92 = @(17.7, 4.1, 0, 77.2, 0.19, 0.8100000000000001, 0.3, 0.49, 0.18)
    93 = @(26.6, 5.2, 0, 67.2, 0.28, 0.84, 0.42, 0.6, 0.27)
    94 = @(9.9, 1.9, 0, 87.2, 0.1, 0.84, 0.18, 0.34, 0.1)
    95 = @(30.9, 5.5, 0, 62.6, 0.33, 0.85, 0.48, 0.65, 0.31)
    96 = @(42.4, 6.3, 0, 50.3, 0.46, 0.87, 0.6, 0.74, 0.43)
    97 = @(15.8, 2.4, 0, 80.8, 0.16, 0.87, 0.28, 0.47, 0.16)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column D is index 4 (1-based), so offset by 4
        $col = $i + 4
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
